# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that were regenerated when the
# handback report was produced again (times shifted by ~1 minute).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file row.
$wsOverview.Range("G2").Value = "2016-09-07 01:25:16"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsZhCn.Range("H2").Value = "2016-09-07 01:25:10"
$wsZhCn.Range("K2").Value = "2016-09-07 01:25:39"

# de-de sheet: Correspond Handoff Datetime (same value as Overview G2) and
# Correspond Handback DateTime.
$wsDeDe.Range("H2").Value = "2016-09-07 01:25:16"
$wsDeDe.Range("K2").Value = "2016-09-07 01:25:47"
